# Rename entries in the systemtools list:
#   "sql"  -> "intellij"
#   "kasp" -> "helpdesk"
# Also move the active selection to D6 (as recorded by Excel on save).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "helpdesk"
$ws.Range("A3").Value = "intellij"

$ws.Range("D6").Select()
